# Update sheet1 data rows (A2:T17) with new TPM-derived values.
# The edit adds a new "Resolving-Mac" sending-cluster block (rows 14-17) and
# refreshes the receptor/edge-weight statistics (columns G-T) for every existing row,
# since those now reflect the larger 4-cluster population.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
$data[0,0] = "ECs"
$data[0,1] = "Wnt16"
$data[0,2] = "Lrp5"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.4491763333333333
$data[0,7] = 1.347529
$data[0,8] = 0.236099761434867
$data[0,9] = 0.236099761434867
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 18.42392
$data[0,13] = 55.27176
$data[0,14] = 0.3903243738016154
$data[0,15] = 0.3903243738016154
$data[0,16] = 8.275588831226667
$data[0,17] = 74.48029948103999
$data[0,18] = 0.09215549153677523
$data[0,19] = 0.09215549153677524
$data[1,0] = "ECs"
$data[1,1] = "Wnt16"
$data[1,2] = "Lrp5"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.4491763333333333
$data[1,7] = 1.347529
$data[1,8] = 0.236099761434867
$data[1,9] = 0.236099761434867
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 12.393653
$data[1,13] = 37.180959
$data[1,14] = 0.2625687066780312
$data[1,15] = 0.2625687066780312
$data[1,16] = 5.566935611145667
$data[1,17] = 50.102420500311
$data[1,18] = 0.06199240900694472
$data[1,19] = 0.06199240900694472
$data[2,0] = "ECs"
$data[2,1] = "Wnt16"
$data[2,2] = "Lrp5"
$data[2,3] = "MuSCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.4491763333333333
$data[2,7] = 1.347529
$data[2,8] = 0.236099761434867
$data[2,9] = 0.236099761434867
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 10.76369066666667
$data[2,13] = 32.291072
$data[2,14] = 0.2280367489253622
$data[2,15] = 0.2280367489253622
$data[2,16] = 4.834795106787556
$data[2,17] = 43.513155961088
$data[2,18] = 0.05383942201966067
$data[2,19] = 0.05383942201966067
$data[3,0] = "ECs"
$data[3,1] = "Wnt16"
$data[3,2] = "Lrp5"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 2
$data[3,5] = 0.6666666666666666
$data[3,6] = 0.4491763333333333
$data[3,7] = 1.347529
$data[3,8] = 0.236099761434867
$data[3,9] = 0.236099761434867
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 5.620297999999999
$data[3,13] = 16.860894
$data[3,14] = 0.1190701705949913
$data[3,15] = 0.1190701705949913
$data[3,16] = 2.524504847880666
$data[3,17] = 22.720543630926
$data[3,18] = 0.02811243887148635
$data[3,19] = 0.02811243887148635
$data[4,0] = "FAPs"
$data[4,1] = "Wnt16"
$data[4,2] = "Lrp5"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.190874
$data[4,7] = 3.572622
$data[4,8] = 0.6259569937989885
$data[4,9] = 0.6259569937989886
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 18.42392
$data[4,13] = 55.27176
$data[4,14] = 0.3903243738016154
$data[4,15] = 0.3903243738016154
$data[4,16] = 21.94056730608
$data[4,17] = 197.46510575472
$data[4,18] = 0.2443262716313318
$data[4,19] = 0.2443262716313319
$data[5,0] = "FAPs"
$data[5,1] = "Wnt16"
$data[5,2] = "Lrp5"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.190874
$data[5,7] = 3.572622
$data[5,8] = 0.6259569937989885
$data[5,9] = 0.6259569937989886
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 12.393653
$data[5,13] = 37.180959
$data[5,14] = 0.2625687066780312
$data[5,15] = 0.2625687066780312
$data[5,16] = 14.759279122722
$data[5,17] = 132.833512104498
$data[5,18] = 0.1643567182978688
$data[5,19] = 0.1643567182978688
$data[6,0] = "FAPs"
$data[6,1] = "Wnt16"
$data[6,2] = "Lrp5"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.190874
$data[6,7] = 3.572622
$data[6,8] = 0.6259569937989885
$data[6,9] = 0.6259569937989886
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 10.76369066666667
$data[6,13] = 32.291072
$data[6,14] = 0.2280367489253622
$data[6,15] = 0.2280367489253622
$data[6,16] = 12.818199358976
$data[6,17] = 115.363794230784
$data[6,18] = 0.1427411978330144
$data[6,19] = 0.1427411978330144
$data[7,0] = "FAPs"
$data[7,1] = "Wnt16"
$data[7,2] = "Lrp5"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.190874
$data[7,7] = 3.572622
$data[7,8] = 0.6259569937989885
$data[7,9] = 0.6259569937989886
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 5.620297999999999
$data[7,13] = 16.860894
$data[7,14] = 0.1190701705949913
$data[7,15] = 0.1190701705949913
$data[7,16] = 6.693066760451999
$data[7,17] = 60.23760084406799
$data[7,18] = 0.07453280603677344
$data[7,19] = 0.07453280603677345
$data[8,0] = "MuSCs"
$data[8,1] = "Wnt16"
$data[8,2] = "Lrp5"
$data[8,3] = "ECs"
$data[8,4] = 1
$data[8,5] = 0.3333333333333333
$data[8,6] = 0.06510733333333334
$data[8,7] = 0.195322
$data[8,8] = 0.03422225243611164
$data[8,9] = 0.03422225243611164
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 18.42392
$data[8,13] = 55.27176
$data[8,14] = 0.3903243738016154
$data[8,15] = 0.3903243738016154
$data[8,16] = 1.199532300746667
$data[8,17] = 10.79579070672
$data[8,18] = 0.01335777925220608
$data[8,19] = 0.01335777925220608
$data[9,0] = "MuSCs"
$data[9,1] = "Wnt16"
$data[9,2] = "Lrp5"
$data[9,3] = "FAPs"
$data[9,4] = 1
$data[9,5] = 0.3333333333333333
$data[9,6] = 0.06510733333333334
$data[9,7] = 0.195322
$data[9,8] = 0.03422225243611164
$data[9,9] = 0.03422225243611164
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 12.393653
$data[9,13] = 37.180959
$data[9,14] = 0.2625687066780312
$data[9,15] = 0.2625687066780312
$data[9,16] = 0.8069176970886668
$data[9,17] = 7.262259273798001
$data[9,18] = 0.008985692561758934
$data[9,19] = 0.008985692561758934
$data[10,0] = "MuSCs"
$data[10,1] = "Wnt16"
$data[10,2] = "Lrp5"
$data[10,3] = "MuSCs"
$data[10,4] = 1
$data[10,5] = 0.3333333333333333
$data[10,6] = 0.06510733333333334
$data[10,7] = 0.195322
$data[10,8] = 0.03422225243611164
$data[10,9] = 0.03422225243611164
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 10.76369066666667
$data[10,13] = 32.291072
$data[10,14] = 0.2280367489253622
$data[10,15] = 0.2280367489253622
$data[10,16] = 0.7007951961315556
$data[10,17] = 6.307156765184
$data[10,18] = 0.007803931186433956
$data[10,19] = 0.007803931186433955
$data[11,0] = "MuSCs"
$data[11,1] = "Wnt16"
$data[11,2] = "Lrp5"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 1
$data[11,5] = 0.3333333333333333
$data[11,6] = 0.06510733333333334
$data[11,7] = 0.195322
$data[11,8] = 0.03422225243611164
$data[11,9] = 0.03422225243611164
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 5.620297999999999
$data[11,13] = 16.860894
$data[11,14] = 0.1190701705949913
$data[11,15] = 0.1190701705949913
$data[11,16] = 0.3659226153186667
$data[11,17] = 3.293303537868
$data[11,18] = 0.004074849435712668
$data[11,19] = 0.004074849435712668
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Wnt16"
$data[12,2] = "Lrp5"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 0.6666666666666666
$data[12,6] = 0.1973276666666667
$data[12,7] = 0.5919829999999999
$data[12,8] = 0.1037209923300328
$data[12,9] = 0.1037209923300328
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 18.42392
$data[12,13] = 55.27176
$data[12,14] = 0.3903243738016154
$data[12,15] = 0.3903243738016154
$data[12,16] = 3.635549144453333
$data[12,17] = 32.71994230008
$data[12,18] = 0.04048483138130222
$data[12,19] = 0.04048483138130223
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Wnt16"
$data[13,2] = "Lrp5"
$data[13,3] = "FAPs"
$data[13,4] = 2
$data[13,5] = 0.6666666666666666
$data[13,6] = 0.1973276666666667
$data[13,7] = 0.5919829999999999
$data[13,8] = 0.1037209923300328
$data[13,9] = 0.1037209923300328
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 12.393653
$data[13,13] = 37.180959
$data[13,14] = 0.2625687066780312
$data[13,15] = 0.2625687066780312
$data[13,16] = 2.445610627966333
$data[13,17] = 22.010495651697
$data[13,18] = 0.02723388681145872
$data[13,19] = 0.02723388681145872
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Wnt16"
$data[14,2] = "Lrp5"
$data[14,3] = "MuSCs"
$data[14,4] = 2
$data[14,5] = 0.6666666666666666
$data[14,6] = 0.1973276666666667
$data[14,7] = 0.5919829999999999
$data[14,8] = 0.1037209923300328
$data[14,9] = 0.1037209923300328
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 10.76369066666667
$data[14,13] = 32.291072
$data[14,14] = 0.2280367489253622
$data[14,15] = 0.2280367489253622
$data[14,16] = 2.123973963975111
$data[14,17] = 19.115765675776
$data[14,18] = 0.02365219788625312
$data[14,19] = 0.02365219788625312
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Wnt16"
$data[15,2] = "Lrp5"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 2
$data[15,5] = 0.6666666666666666
$data[15,6] = 0.1973276666666667
$data[15,7] = 0.5919829999999999
$data[15,8] = 0.1037209923300328
$data[15,9] = 0.1037209923300328
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 5.620297999999999
$data[15,13] = 16.860894
$data[15,14] = 0.1190701705949913
$data[15,15] = 0.1190701705949913
$data[15,16] = 1.109040290311333
$data[15,17] = 9.981362612801998
$data[15,18] = 0.01235007625101879
$data[15,19] = 0.01235007625101879

# Rows 2-17 (1-based Excel rows), columns A-T (1-20); existing header row (row 1) is untouched.
$ws.Range("A2:T17").Value = $data

Write-Output "Updated $($ws.Range('A2:T17').Rows.Count) data rows."
